# Raport.xlsx edit: add a new completed work item ("AdresyControllerTests.cs",
# 34h) to the "Testy" block on Arkusz1, and update the sheet's view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Extend the "Testy" merged label column (H10:H19 -> H10:H20) ---------
# Row 20 previously had nothing in columns H/I/J; copy the blank, centered
# date-style formatting from H19 down into the new H20 cell, then re-merge
# the whole H10:H20 block so it keeps spanning the (now taller) section.
$ws.Range("H19").Copy($ws.Range("H20"))
$ws.Range("H10:H20").Merge()

# --- New row of data: AdresyControllerTests.cs, 34 hours ------------------
$ws.Range("I20").Value = "AdresyControllerTests.cs"
$ws.Range("J20").Value = 34

# --- Update the saved view: scroll so column C is leftmost, and move the
# active selection to J21 -------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 3
$ws.Range("J21").Select()

# Recalculate so the SUM()/share formulas in row 3 & 4 pick up the new
# J20 value (J3 546 -> 580, and the dependent D4/G4/J4/M4 ratios).
$excel.CalculateFull()
